$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 13 de Agosto de 2020 a las 20:08'
$ws.Range("B4").Value = 5382741
$ws.Range("C4").Value = 22439
$ws.Range("D4").Value = 2815723
$ws.Range("E4").Value = 2397408
$ws.Range("G4").Value = 479
$ws.Range("H4").Value = 169610
$ws.Range("B6").Value = 2459464
$ws.Range("C6").Value = 63993
$ws.Range("D6").Value = 1750568
$ws.Range("E6").Value = 660754
$ws.Range("G6").Value = 1004
$ws.Range("H6").Value = 48142
$ws.Range("B12").Value = 380034
$ws.Range("C12").Value = 1866
$ws.Range("D12").Value = 353131
$ws.Range("E12").Value = 16604
$ws.Range("G12").Value = 94
$ws.Range("H12").Value = 10299
$ws.Range("B13").Value = 379799
$ws.Range("C13").Value = 2935
$ws.Range("G13").Value = 26
$ws.Range("H13").Value = 28605
$ws.Range("B21").Value = 245635
$ws.Range("C21").Value = 1243
$ws.Range("D21").Value = 228057
$ws.Range("E21").Value = 11666
$ws.Range("G21").Value = 21
$ws.Range("H21").Value = 5912
$ws.Range("B22").Value = 221977
$ws.Range("C22").Value = 1127
$ws.Range("E22").Value = 11896
$ws.Range("B23").Value = 209365
$ws.Range("C23").Value = 2669
$ws.Range("E23").Value = 95505
$ws.Range("G23").Value = 17
$ws.Range("H23").Value = 30388
$ws.Range("B24").Value = 164277
$ws.Range("C24").Value = 3841
$ws.Range("D24").Value = 117208
$ws.Range("E24").Value = 41428
$ws.Range("G24").Value = 53
$ws.Range("H24").Value = 5641
$ws.Range("B28").Value = 114281
$ws.Range("C28").Value = 343
$ws.Range("D28").Value = 110957
$ws.Range("E28").Value = 3134
$ws.Range("B33").Value = 89555
$ws.Range("C33").Value = 1404
$ws.Range("D33").Value = 64721
$ws.Range("E33").Value = 24183
$ws.Range("G33").Value = 12
$ws.Range("H33").Value = 651
$ws.Range("A57").Value = 'Marruecos'
$ws.Range("B57").Value = 37935
$ws.Range("C57").Value = 1241
$ws.Range("D57").Value = 26687
$ws.Range("E57").Value = 10664
$ws.Range("G57").Value = 28
$ws.Range("H57").Value = 584
$ws.Range("A58").Value = 'Afganistan'
$ws.Range("B58").Value = 37424
$ws.Range("C58").Value = 79
$ws.Range("D58").Value = 26714
$ws.Range("E58").Value = 9347
$ws.Range("G58").Value = 9
$ws.Range("H58").Value = 1363
$ws.Range("A59").Value = 'Suiza'
$ws.Range("B59").Value = 37403
$ws.Range("C59").Value = 234
$ws.Range("D59").Value = 32700
$ws.Range("E59").Value = 2712
$ws.Range("H59").Value = 1991
$ws.Range("A60").Value = 'Argelia'
$ws.Range("B60").Value = 36699
$ws.Range("D60").Value = 25627
$ws.Range("E60").Value = 9739
$ws.Range("H60").Value = 1333
$ws.Range("B67").Value = 26929
$ws.Range("C67").Value = 91
$ws.Range("E67").Value = 1791
$ws.Range("B68").Value = 26204
$ws.Range("C68").Value = 1086
$ws.Range("D68").Value = 11428
$ws.Range("E68").Value = 14297
$ws.Range("G68").Value = 16
$ws.Range("H68").Value = 479
$ws.Range("A115").Value = 'Namibia'
$ws.Range("B115").Value = 3544
$ws.Range("C115").Value = 138
$ws.Range("D115").Value = 848
$ws.Range("E115").Value = 2669
$ws.Range("G115").Value = 5
$ws.Range("H115").Value = 27
$ws.Range("A116").Value = 'Suazilandia'
$ws.Range("B116").Value = 3525
$ws.Range("D116").Value = 1910
$ws.Range("E116").Value = 1552
$ws.Range("H116").Value = 63
$ws.Range("B135").Value = 1940
$ws.Range("C135").Value = 3
$ws.Range("D135").Value = 1496
$ws.Range("E135").Value = 375
$ws.Range("B162").Value = 754
$ws.Range("C162").Value = 20
$ws.Range("E162").Value = 118
